$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at row 68 ---------------------------------------
$ws.Rows.Item(68).Insert()

$ws.Cells.Item(68,1).Value  = 10
$ws.Cells.Item(68,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(68,3).Value  = "La Araucanía"
$ws.Cells.Item(68,4).Value  = 44748
$ws.Cells.Item(68,5).Value  = 9
$ws.Cells.Item(68,6).Value  = "Fruta"
$ws.Cells.Item(68,7).Value  = 100108
$ws.Cells.Item(68,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(68,9).Value  = 100108007
$ws.Cells.Item(68,10).Value = "Coco"
$ws.Cells.Item(68,11).Value = "Sin especificar"
$ws.Cells.Item(68,12).Value = "Primera"
$ws.Cells.Item(68,13).Value = 40
$ws.Cells.Item(68,14).Value = 28000
$ws.Cells.Item(68,15).Value = 28000
$ws.Cells.Item(68,16).Value = 28000
$ws.Cells.Item(68,17).Value = "$/malla 20 unidades"
$ws.Cells.Item(68,18).Value = "Perú"
$ws.Cells.Item(68,19).Value = 1400
$ws.Cells.Item(68,20).Value = 20

# --- Insert second new row at row 75 (post first-shift numbering) ---------
$ws.Rows.Item(75).Insert()

$ws.Cells.Item(75,1).Value  = 10
$ws.Cells.Item(75,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(75,3).Value  = "La Araucanía"
$ws.Cells.Item(75,4).Value  = 44747
$ws.Cells.Item(75,5).Value  = 9
$ws.Cells.Item(75,6).Value  = "Fruta"
$ws.Cells.Item(75,7).Value  = 100108
$ws.Cells.Item(75,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(75,9).Value  = 100108007
$ws.Cells.Item(75,10).Value = "Coco"
$ws.Cells.Item(75,11).Value = "Sin especificar"
$ws.Cells.Item(75,12).Value = "Primera"
$ws.Cells.Item(75,13).Value = 20
$ws.Cells.Item(75,14).Value = 28000
$ws.Cells.Item(75,15).Value = 28000
$ws.Cells.Item(75,16).Value = 28000
$ws.Cells.Item(75,17).Value = "$/malla 20 unidades"
$ws.Cells.Item(75,18).Value = "Perú"
$ws.Cells.Item(75,19).Value = 1400
$ws.Cells.Item(75,20).Value = 20

Write-Output "done"
